# Necron_units.xlsx edit script
# Updates wargear "Options" text on several sheets (HQ, Troops, Elites,
# Fast Attack) to fix wording/typos and add a missing option, and moves
# each sheet's active-cell selection to match the saved state.

$wb = $excel.ActiveWorkbook

# --- HQ sheet ---
$ws = $wb.Worksheets.Item("HQ")
$ws.Range("E2").Value = "Gauss cannon/Tesla cannon, Staff of light/Warscythe/Hyperphase sword/Voidblade, Phylactery, Resurrection orb"
$ws.Range("E3").Value = "Phylactery, Canoptek cloak/Chronometron"
$ws.Range("E3").Select()

# --- Troops sheet ---
$ws = $wb.Worksheets.Item("Troops")
$ws.Range("E2").Value = "Gauss blaster/Tesla carbine"
$ws.Range("E2").Select()

# --- Fast Attack sheet ---
$ws = $wb.Worksheets.Item("Fast Attack")
$ws.Range("E4").Value = "Gauss cannon/Heavy gauss cannon-3"
$ws.Range("G12").Select()

# --- Elites sheet (kept as the active tab, matching the saved workbook) ---
$ws = $wb.Worksheets.Item("Elites")
$ws.Range("E5").Value = "Rod of covenant/Particle caster+Voidblade"
$ws.Range("E6").Value = "Heat ray/2*Heavy gauss cannon/Particle shredder"
$ws.Range("E6").Select()
